# Updates the sheet name, reference date, and several Saldo Previsto / Vl. Total
# values to reflect the newer IClientBalance export (2024-08-30 090752).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the newer export timestamp.
$ws.Name = "IClientBalance-20240830-090752-"

# Column G ("Dt. Referencia") moves from 45532 (2024-08-28) to 45534 (2024-08-30)
# for every data row (rows 2 through 274).
$lastRow = 274
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 45534
}

# Rows whose "Saldo Previsto" (E) / "Vl. Total" (H) values changed between the
# two export runs. Both columns always carry the same value for a given row.
$valueChanges = @{
    5 = 1118.97
    8 = 10058.83
    15 = 81547.55
    17 = 1153.5
    43 = 1269.5
    49 = 7076.3
    52 = 148557.06
    57 = 49288.66
    58 = 8769.13
    60 = 701.65
    99 = 1355.22
    102 = 7987.23
    103 = 999.99
    104 = 2255.01
    108 = 2312.35
    112 = 799.71
    120 = 0
    132 = 1050.85
    143 = 2365.78
    158 = 112.94
    161 = 271.55
    173 = 1700.44
    189 = 19196.12
    224 = 541.29
    230 = 55939.7
    231 = 821.06
    232 = 218.17
    235 = 350.71
    249 = 974.54
    255 = 16680.47
    264 = 1986.68
    265 = 1444.82
    270 = -229.38
    271 = 922.52
    273 = 1291.08
}

foreach ($row in $valueChanges.Keys) {
    $val = $valueChanges[$row]
    $ws.Cells.Item($row, 5).Value = $val
    $ws.Cells.Item($row, 8).Value = $val
}
